$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "27.455.52"
$ws.Range("E2").Value = "  -1.48%  "
Set-TextValue "D3" "1.746.75"
$ws.Range("E3").Value = "  -1.39%  "
Set-TextValue "D4" "1.003"
$ws.Range("E4").Value = "  +0.56%  "
Set-TextValue "D5" "324.98"
$ws.Range("E5").Value = "  +0.90%  "
Set-TextValue "D7" "0.4444"
$ws.Range("E7").Value = "  +4.25%  "
Set-TextValue "D8" "0.3566"
$ws.Range("E8").Value = "  -1.21%  "
Set-TextValue "D9" "0.07474"
$ws.Range("E9").Value = "  +0.16%  "
Set-TextValue "D10" "41.98"
$ws.Range("E10").Value = "  -5.29%  "
Set-TextValue "D11" "1.087"
$ws.Range("E11").Value = "  -1.54%  "
Set-TextValue "D12" "1.002"
$ws.Range("E12").Value = "  +0.78%  "
Set-TextValue "D13" "20.65"
$ws.Range("E13").Value = "  -4.17%  "
Set-TextValue "D14" "5.998"
$ws.Range("E14").Value = "  -2.04%  "
Set-TextValue "D15" "7.080"
$ws.Range("E15").Value = "  -3.04%  "
Set-TextValue "D16" "1.752.55"
Set-TextValue "D17" "92.70"
$ws.Range("E17").Value = "  +1.67%  "
Set-TextValue "D18" "0.00001057"
$ws.Range("E18").Value = "  -0.28%  "
Set-TextValue "D19" "0.06408"
$ws.Range("E19").Value = "  +0.92%  "
$ws.Range("E20").Value = "  +0.50%  "
Set-TextValue "D21" "16.76"
$ws.Range("E21").Value = "  -2.36%  "
Set-TextValue "D22" "5.797"
$ws.Range("E22").Value = "  -2.42%  "
Set-TextValue "D23" "27.519.33"
$ws.Range("E23").Value = "  -1.30%  "
Set-TextValue "D24" "11.14"
$ws.Range("E24").Value = "  -1.88%  "
Set-TextValue "D25" "2.096"
$ws.Range("E25").Value = "  -3.07%  "
Set-TextValue "D26" "163.15"
$ws.Range("E26").Value = "  +1.97%  "
Set-TextValue "D27" "20.39"
$ws.Range("E27").Value = "  +0.67%  "
Set-TextValue "D28" "1.949.41"
$ws.Range("E28").Value = "  -2.11%  "
Set-TextValue "D29" "2.065"
$ws.Range("E29").Value = "  -4.58%  "
Set-TextValue "D30" "125.37"
$ws.Range("E30").Value = "  -0.30%  "
Set-TextValue "D31" "1.069"
$ws.Range("E31").Value = "  -8.52%  "
Set-TextValue "D32" "3.672"
$ws.Range("E32").Value = "  +4.28%  "
Set-TextValue "D33" "0.09044"
$ws.Range("E33").Value = "  +0.63%  "
Set-TextValue "D34" "5.479"
$ws.Range("E34").Value = "  -3.60%  "
Set-TextValue "D35" "11.89"
$ws.Range("E35").Value = "  -5.99%  "
Set-TextValue "D36" "0.02279"
$ws.Range("E36").Value = "  -1.54%  "
Set-TextValue "D37" "0.2090"
$ws.Range("E37").Value = "  -0.93%  "
Set-TextValue "D38" "0.6344"
$ws.Range("E38").Value = "  -1.26%  "
Set-TextValue "D39" "0.05996"
$ws.Range("E39").Value = "  -0.88%  "
Set-TextValue "D40" "4.921"
$ws.Range("E40").Value = "  -2.91%  "
Set-TextValue "D41" "1.201"
$ws.Range("E41").Value = "  +1.94%  "
Set-TextValue "D42" "1.378"
$ws.Range("E42").Value = "  -0.72%  "
Set-TextValue "D43" "7.714"
$ws.Range("E43").Value = "  -1.76%  "
Set-TextValue "D44" "13.23"
$ws.Range("E44").Value = "  -3.02%  "
$ws.Range("E45").Value = "  +0.59%  "
Set-TextValue "D46" "0.5874"
$ws.Range("E46").Value = "  -1.48%  "
Set-TextValue "D47" "121.42"
$ws.Range("E47").Value = "  -2.17%  "
Set-TextValue "D48" "1.942"
$ws.Range("E48").Value = "  -1.98%  "
Set-TextValue "D49" "1.139"
$ws.Range("E49").Value = "  -0.65%  "
Set-TextValue "D50" "0.06825"
$ws.Range("E50").Value = "  -0.82%  "
Set-TextValue "D51" "72.09"
$ws.Range("E51").Value = "  -3.17%  "
